$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ají" (Vega Monumental Concepción).
# Insert a new row at row 41 so the existing rows 41-46 shift down to 42-47,
# then populate the new row 41 with the latest record's data.
$ws.Rows.Item(41).Insert()

$ws.Cells.Item(41, 1).Value = 11
$ws.Cells.Item(41, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(41, 3).Value = "Bíobío"
$ws.Cells.Item(41, 4).Value = 44474
$ws.Cells.Item(41, 5).Value = 8
$ws.Cells.Item(41, 6).Value = 100112021
$ws.Cells.Item(41, 7).Value = "Ají"
$ws.Cells.Item(41, 8).Value = "Inferno"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 40
$ws.Cells.Item(41, 11).Value = 40000
$ws.Cells.Item(41, 12).Value = 42000
$ws.Cells.Item(41, 13).Value = 41000
$ws.Cells.Item(41, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(41, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(41, 16).Value = 3417
$ws.Cells.Item(41, 17).Value = 12
$ws.Cells.Item(41, 18).Value = "Hortaliza"
